# Insert a new "toilet location" variable row into the vulnerabilities
# data dictionary (row 309), which pushes the existing rows 309:501 down
# to 310:502, then fix up the bookkeeping (AutoFilter range, the
# _FilterDatabase defined name, and the view/selection) so the workbook
# matches where the author ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vulnerabilities")

# Insert a new row above current row 309; this shifts rows 309-501 down to
# 310-502 and keeps all their values/formatting intact.
$ws.Rows.Item(309).Insert() | Out-Null

# Populate the newly inserted row 309 with the new "toilet location"
# variable entry (matches the author's addition to the data dictionary).
$ws.Range("A309").Value = "toilet.loc"
$ws.Range("B309").Value = "Toilet Location"
$ws.Range("C309").Value = "Location of the household toilet"
$ws.Range("F309").Value = 1

# The sheet's AutoFilter still covers only the old A1:I501 range after the
# insert, so turn it off and reapply it across the new full extent
# (A1:I502) to pick up the newly-added row.
$ws.AutoFilterMode = $false
$ws.Range("A1:I502").AutoFilter() | Out-Null

# Keep the workbook-level _FilterDatabase defined name (used by the
# AutoFilter UI state) in sync with the new range too.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "vulnerabilities!_FilterDatabase") {
        $n.RefersTo = "=vulnerabilities!`$A`$1:`$I`$502"
    }
}

# Activate the sheet and select the cell the author ended up on (just
# below / right of the newly-inserted row) so the saved selection matches.
$ws.Activate()
$ws.Range("F310").Select()
